$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.790.46'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").Value = '1.701.37'
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3954'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4083'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.60%  '

$ws.Range("E9").Value = '  -1.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.86'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08914'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.697'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.145'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("E16").Value = '  +1.17%  '

$ws.Range("D17").Value = '1.711.35'
$ws.Range("E17").Value = '  +1.08%  '

$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07142'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.02%  '

$ws.Range("E20").Value = '  +2.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.224'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.57%  '

$ws.Range("D24").Value = '24.783.46'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.096'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.339'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.339'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +24.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '164.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.197'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.068'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09183'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.081'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03076'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2812'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.69%  '

$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.966'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.88%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.54%  '

$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09307'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7836'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.476'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.641'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7265'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.246'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.356'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.15%  '

$ws.Range("E48").Value = '  +0.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '93.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08067'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
